# chore: update Sheets via scheduled runner
# Refreshes the cached market-price / leve-profit figures (columns H-N) on a
# handful of rows across each crafting-job sheet (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR) of the Cactuar_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 859.9792
$ws.Range("J17").Value = 859.9792
$ws.Range("L17").Value = 2579.9376
$ws.Range("N17").Value = -2915.9376

$ws.Range("H51").Value = 5062.8
$ws.Range("J51").Value = 4832.4443
$ws.Range("L51").Value = 4832.4443
$ws.Range("N51").Value = -5800.4443

$ws.Range("H62").Value = 3844.3333
$ws.Range("I62").Value = 3525
$ws.Range("J62").Value = 4099.8
$ws.Range("K62").Value = 3525
$ws.Range("L62").Value = 4099.8
$ws.Range("M62").Value = -2901
$ws.Range("N62").Value = -5347.8

$ws.Range("H65").Value = 3844.3333
$ws.Range("I65").Value = 3525
$ws.Range("J65").Value = 4099.8
$ws.Range("K65").Value = 17625
$ws.Range("L65").Value = 20499
$ws.Range("M65").Value = -14505
$ws.Range("N65").Value = -26739

$ws.Range("H100").Value = 1694.5
$ws.Range("I100").Value = 1694.5
$ws.Range("K100").Value = 1694.5
$ws.Range("M100").Value = -1153.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 648593.5600000001
$ws.Range("I2").Value = 1029306.9
$ws.Range("J2").Value = 1381
$ws.Range("K2").Value = 1029306.9
$ws.Range("L2").Value = 1381
$ws.Range("M2").Value = -1029193.9
$ws.Range("N2").Value = -1607

$ws.Range("H32").Value = 14207.955
$ws.Range("I32").Value = 14703.5
$ws.Range("K32").Value = 14703.5
$ws.Range("M32").Value = -14416.5

$ws.Range("H61").Value = 8384.125
$ws.Range("I61").Value = 9588.191999999999
$ws.Range("K61").Value = 9588.191999999999
$ws.Range("M61").Value = -9376.191999999999

$ws.Range("H116").Value = 648593.5600000001
$ws.Range("I116").Value = 1029306.9
$ws.Range("J116").Value = 1381
$ws.Range("K116").Value = 1029306.9
$ws.Range("L116").Value = 1381
$ws.Range("M116").Value = -1027012.9
$ws.Range("N116").Value = -5969

$ws.Range("H132").Value = 13603.1
$ws.Range("I132").Value = 15839.872
$ws.Range("K132").Value = 47519.61599999999
$ws.Range("M132").Value = -44989.61599999999

$ws.Range("H136").Value = 8384.125
$ws.Range("I136").Value = 9588.191999999999
$ws.Range("K136").Value = 28764.576
$ws.Range("M136").Value = -26214.576

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 648593.5600000001
$ws.Range("I3").Value = 1029306.9
$ws.Range("J3").Value = 1381
$ws.Range("K3").Value = 1029306.9
$ws.Range("L3").Value = 1381
$ws.Range("M3").Value = -1029192.9
$ws.Range("N3").Value = -1609

$ws.Range("H105").Value = 1957.1364
$ws.Range("I105").Value = 1901.3334
$ws.Range("J105").Value = 2076.7144
$ws.Range("K105").Value = 1901.3334
$ws.Range("L105").Value = 2076.7144
$ws.Range("M105").Value = -154.3334
$ws.Range("N105").Value = -5570.7144

$ws.Range("H107").Value = 3241.7
$ws.Range("I107").Value = 3554.625
$ws.Range("K107").Value = 3554.625
$ws.Range("M107").Value = -1634.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 282.7143
$ws.Range("I7").Value = 220
$ws.Range("K7").Value = 220
$ws.Range("M7").Value = -107

$ws.Range("H86").Value = 7605.45
$ws.Range("J86").Value = 9205.429
$ws.Range("L86").Value = 9205.429
$ws.Range("N86").Value = -11451.429

$ws.Range("H89").Value = 7605.45
$ws.Range("J89").Value = 9205.429
$ws.Range("L89").Value = 46027.145
$ws.Range("N89").Value = -57259.145

$ws.Range("H99").Value = 11072.728
$ws.Range("J99").Value = 11427.857
$ws.Range("L99").Value = 11427.857
$ws.Range("N99").Value = -14423.857

$ws.Range("H122").Value = 5106.5386
$ws.Range("I122").Value = 2055
$ws.Range("K122").Value = 6165
$ws.Range("M122").Value = -3715

$ws.Range("H126").Value = 11072.728
$ws.Range("J126").Value = 11427.857
$ws.Range("L126").Value = 34283.571
$ws.Range("N126").Value = -39223.571

$ws.Range("H134").Value = 2456.913
$ws.Range("I134").Value = 1939.4375
$ws.Range("K134").Value = 5818.3125
$ws.Range("M134").Value = -3283.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5735552.5
$ws.Range("I4").Value = 1454844.6
$ws.Range("K4").Value = 4364533.800000001
$ws.Range("M4").Value = -4364421.800000001

$ws.Range("H7").Value = 287.45456
$ws.Range("J7").Value = 233
$ws.Range("L7").Value = 699
$ws.Range("N7").Value = -923

$ws.Range("H9").Value = 476.66666
$ws.Range("J9").Value = 665.5
$ws.Range("L9").Value = 1996.5
$ws.Range("N9").Value = -2444.5

$ws.Range("H21").Value = 383.66666
$ws.Range("I21").Value = 383.66666
$ws.Range("K21").Value = 1150.99998
$ws.Range("M21").Value = -977.9999800000001

$ws.Range("H34").Value = 689640.4399999999
$ws.Range("I34").Value = 1033280.7
$ws.Range("J34").Value = 2360
$ws.Range("K34").Value = 3099842.1
$ws.Range("L34").Value = 7080
$ws.Range("M34").Value = -3099758.1
$ws.Range("N34").Value = -7248

$ws.Range("H46").Value = 2285.7144
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2909

$ws.Range("H107").Value = 2776.6667
$ws.Range("J107").Value = 1532
$ws.Range("L107").Value = 4596
$ws.Range("N107").Value = -8436

$ws.Range("H122").Value = 859.8
$ws.Range("I122").Value = 699.75
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6297.75
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -3847.75
$ws.Range("N122").Value = -18400

$ws.Range("H137").Value = 53584356
$ws.Range("I137").Value = 75001896
$ws.Range("K137").Value = 225005688
$ws.Range("M137").Value = -225000588

$ws.Range("H140").Value = 8422.947
$ws.Range("I140").Value = 3639.7273
$ws.Range("K140").Value = 10919.1819
$ws.Range("M140").Value = -5739.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3035779.2
$ws.Range("I70").Value = 3501722.2
$ws.Range("K70").Value = 3501722.2
$ws.Range("M70").Value = -3501452.2

$ws.Range("H73").Value = 3035779.2
$ws.Range("I73").Value = 3501722.2
$ws.Range("K73").Value = 3501722.2
$ws.Range("M73").Value = -3500786.2

$ws.Range("H123").Value = 46107.152
$ws.Range("J123").Value = 46107.152
$ws.Range("L123").Value = 46107.152
$ws.Range("N123").Value = -51007.152

$ws.Range("H132").Value = 41586.37
$ws.Range("I132").Value = 52619.242
$ws.Range("J132").Value = 6790.385
$ws.Range("K132").Value = 157857.726
$ws.Range("L132").Value = 20371.155
$ws.Range("M132").Value = -155327.726
$ws.Range("N132").Value = -25431.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3628.3914
$ws.Range("I93").Value = 918.6842
$ws.Range("J93").Value = 16499.5
$ws.Range("K93").Value = 918.6842
$ws.Range("L93").Value = 16499.5
$ws.Range("M93").Value = 329.3158
$ws.Range("N93").Value = -18995.5

$ws.Range("H132").Value = 5079.5864
$ws.Range("I132").Value = 4251.6113
$ws.Range("J132").Value = 6434.4546
$ws.Range("K132").Value = 12754.8339
$ws.Range("L132").Value = 19303.3638
$ws.Range("M132").Value = -10224.8339
$ws.Range("N132").Value = -24363.3638

$ws.Range("H136").Value = 2830.4827
$ws.Range("I136").Value = 2182.4167
$ws.Range("K136").Value = 6547.250100000001
$ws.Range("M136").Value = -3997.250100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2781.3
$ws.Range("I107").Value = 2723
$ws.Range("J107").Value = 2868.75
$ws.Range("K107").Value = 8169
$ws.Range("L107").Value = 8606.25
$ws.Range("M107").Value = -6249
$ws.Range("N107").Value = -12446.25

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 1237702
$ws.Range("I132").Value = 2315592.8
$ws.Range("K132").Value = 6946778.399999999
$ws.Range("M132").Value = -6944248.399999999

$ws.Range("H136").Value = 7390.891
$ws.Range("I136").Value = 1875.3939
$ws.Range("K136").Value = 5626.1817
$ws.Range("M136").Value = -3076.1817
